# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data refresh values to the Sheets workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 4556.8
$ws.Range("I31").Value = 3199
$ws.Range("K31").Value = 9597
$ws.Range("M31").Value = -9367

$ws.Range("H86").Value = 66111.375
$ws.Range("I86").Value = 86648.836
$ws.Range("J86").Value = 4499
$ws.Range("K86").Value = 86648.836
$ws.Range("L86").Value = 4499
$ws.Range("M86").Value = -85525.836
$ws.Range("N86").Value = -6745

$ws.Range("H89").Value = 66111.375
$ws.Range("I89").Value = 86648.836
$ws.Range("J89").Value = 4499
$ws.Range("K89").Value = 433244.18
$ws.Range("L89").Value = 22495
$ws.Range("M89").Value = -427628.18
$ws.Range("N89").Value = -33727

$ws.Range("H116").Value = 24641954
$ws.Range("I116").Value = 43594890
$ws.Range("J116").Value = 3139.1
$ws.Range("K116").Value = 43594890
$ws.Range("L116").Value = 3139.1
$ws.Range("M116").Value = -43591448
$ws.Range("N116").Value = -10023.1

$ws.Range("H129").Value = 2539.7036
$ws.Range("J129").Value = 3509.7646
$ws.Range("L129").Value = 10529.2938
$ws.Range("N129").Value = -20529.2938

$ws.Range("H132").Value = 6195.1665
$ws.Range("I132").Value = 1256.7059
$ws.Range("K132").Value = 3770.1177
$ws.Range("M132").Value = -1240.1177

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15771.819
$ws.Range("I32").Value = 15246.887
$ws.Range("K32").Value = 15246.887
$ws.Range("M32").Value = -14959.887

$ws.Range("H61").Value = 4449.1177
$ws.Range("I61").Value = 3551.3142
$ws.Range("K61").Value = 3551.3142
$ws.Range("M61").Value = -3339.3142

$ws.Range("H122").Value = 6120.095
$ws.Range("I122").Value = 4688.9375
$ws.Range("K122").Value = 14066.8125
$ws.Range("M122").Value = -11616.8125

$ws.Range("H132").Value = 17865.162
$ws.Range("I132").Value = 23210.46
$ws.Range("J132").Value = 9690
$ws.Range("K132").Value = 69631.38
$ws.Range("L132").Value = 29070
$ws.Range("M132").Value = -67101.38
$ws.Range("N132").Value = -34130

$ws.Range("H133").Value = 75050.8
$ws.Range("J133").Value = 75050.8
$ws.Range("L133").Value = 75050.8
$ws.Range("N133").Value = -80110.8

$ws.Range("H135").Value = 71742.75
$ws.Range("J135").Value = 71742.75
$ws.Range("L135").Value = 71742.75
$ws.Range("N135").Value = -81882.75

$ws.Range("H136").Value = 4449.1177
$ws.Range("I136").Value = 3551.3142
$ws.Range("K136").Value = 10653.9426
$ws.Range("M136").Value = -8103.942599999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1703.3871
$ws.Range("I134").Value = 1303.138
$ws.Range("J134").Value = 7507
$ws.Range("K134").Value = 3909.414
$ws.Range("L134").Value = 22521
$ws.Range("M134").Value = -1374.414
$ws.Range("N134").Value = -27591

$ws.Range("H135").Value = 122447.5
$ws.Range("J135").Value = 122447.5
$ws.Range("L135").Value = 122447.5
$ws.Range("N135").Value = -132587.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 58999
$ws.Range("J18").Value = 58999
$ws.Range("L18").Value = 58999
$ws.Range("N18").Value = -59459

$ws.Range("H31").Value = 1898.3695
$ws.Range("I31").Value = 1572.381
$ws.Range("K31").Value = 1572.381
$ws.Range("M31").Value = -1277.381

$ws.Range("H34").Value = 1898.3695
$ws.Range("I34").Value = 1572.381
$ws.Range("K34").Value = 1572.381
$ws.Range("M34").Value = -1370.381

$ws.Range("H58").Value = 1127.4546
$ws.Range("I58").Value = 1155.7778
$ws.Range("J58").Value = 1000
$ws.Range("K58").Value = 1155.7778
$ws.Range("L58").Value = 1000
$ws.Range("M58").Value = -952.7778000000001
$ws.Range("N58").Value = -1406

$ws.Range("H105").Value = 1581.1666
$ws.Range("I105").Value = 1559.875
$ws.Range("J105").Value = 1623.75
$ws.Range("K105").Value = 1559.875
$ws.Range("L105").Value = 1623.75
$ws.Range("M105").Value = 187.125
$ws.Range("N105").Value = -5117.75

$ws.Range("H107").Value = 1041.8695
$ws.Range("I107").Value = 1027.6666
$ws.Range("J107").Value = 1051
$ws.Range("K107").Value = 1027.6666
$ws.Range("L107").Value = 1051
$ws.Range("M107").Value = 892.3334
$ws.Range("N107").Value = -4891

$ws.Range("H132").Value = 13339539
$ws.Range("I132").Value = 20835292
$ws.Range("K132").Value = 62505876
$ws.Range("M132").Value = -62503346

$ws.Range("H134").Value = 4374.9375
$ws.Range("I134").Value = 4652.25
$ws.Range("J134").Value = 4282.5
$ws.Range("K134").Value = 13956.75
$ws.Range("L134").Value = 12847.5
$ws.Range("M134").Value = -11421.75
$ws.Range("N134").Value = -17917.5

$ws.Range("H136").Value = 1127.4546
$ws.Range("I136").Value = 1155.7778
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 3467.3334
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -917.3334000000004
$ws.Range("N136").Value = -8100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 118379.4
$ws.Range("I128").Value = 118379.4
$ws.Range("K128").Value = 355138.2
$ws.Range("M128").Value = -350158.2

$ws.Range("H129").Value = 1440.4546
$ws.Range("I129").Value = 964.3333
$ws.Range("J129").Value = 2460.7144
$ws.Range("K129").Value = 2892.9999
$ws.Range("L129").Value = 7382.1432
$ws.Range("M129").Value = 2107.0001
$ws.Range("N129").Value = -17382.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 22022.223
$ws.Range("J15").Value = 21650
$ws.Range("L15").Value = 21650
$ws.Range("N15").Value = -22226

$ws.Range("H81").Value = 22022.223
$ws.Range("J81").Value = 21650
$ws.Range("L81").Value = 21650
$ws.Range("N81").Value = -23646

$ws.Range("H84").Value = 22022.223
$ws.Range("J84").Value = 21650
$ws.Range("L84").Value = 64950
$ws.Range("N84").Value = -74934

$ws.Range("H102").Value = 18525750
$ws.Range("I102").Value = 23816428
$ws.Range("K102").Value = 23816428
$ws.Range("M102").Value = -23814806

$ws.Range("H122").Value = 389616.94
$ws.Range("I122").Value = 717003.4
$ws.Range("K122").Value = 2151010.2
$ws.Range("M122").Value = -2148560.2

$ws.Range("H132").Value = 6996.6665
$ws.Range("I132").Value = 4990
$ws.Range("K132").Value = 14970
$ws.Range("M132").Value = -12440

$ws.Range("H134").Value = 600403
$ws.Range("J134").Value = 600403
$ws.Range("L134").Value = 1801209
$ws.Range("N134").Value = -1806279

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H54").Value = 39303.5
$ws.Range("J54").Value = 39303.5
$ws.Range("L54").Value = 39303.5
$ws.Range("N54").Value = -40591.5

$ws.Range("H61").Value = 1620.3334
$ws.Range("I61").Value = 1529.4667
$ws.Range("K61").Value = 1529.4667
$ws.Range("M61").Value = -1327.4667

$ws.Range("H106").Value = 16499
$ws.Range("J106").Value = 16499
$ws.Range("L106").Value = 16499
$ws.Range("N106").Value = -19023

$ws.Range("H113").Value = 1620.3334
$ws.Range("I113").Value = 1529.4667
$ws.Range("K113").Value = 1529.4667
$ws.Range("M113").Value = 640.5333000000001

$ws.Range("H132").Value = 2658.182
$ws.Range("I132").Value = 2612.3206
$ws.Range("J132").Value = 2828.524
$ws.Range("K132").Value = 7836.9618
$ws.Range("L132").Value = 8485.572
$ws.Range("M132").Value = -5306.9618
$ws.Range("N132").Value = -13545.572

$ws.Range("H133").Value = 92326
$ws.Range("J133").Value = 92326
$ws.Range("L133").Value = 92326
$ws.Range("N133").Value = -97386

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1938.5385
$ws.Range("J132").Value = 2227.6428
$ws.Range("L132").Value = 6682.928400000001
$ws.Range("N132").Value = -11742.9284

$ws.Range("H133").Value = 91333
$ws.Range("J133").Value = 91333
$ws.Range("L133").Value = 91333
$ws.Range("N133").Value = -101453

$ws.Range("H140").Value = 89664.25
$ws.Range("J140").Value = 89664.25
$ws.Range("L140").Value = 89664.25
$ws.Range("N140").Value = -100024.25
